$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder the faculty list and expand it with the full set of faculties/schools ---

# Row 2 (facultyID 1): was SoC, now becomes FASS (reuses existing strings)
$ws.Range("B2").Value = "FASS"
$ws.Range("C2").Value = "Faculty of Arts and Social Sciences"

# Row 3 (facultyID 2): was FASS, now becomes Business (new strings)
$ws.Range("B3").Value = "Business"
$ws.Range("C3").Value = "Business School"

# Row 4 (facultyID 3): was Others/Others, now becomes SoC (reuses existing strings,
# and frees up the old "Others" entry since both its references are rewritten here)
$ws.Range("B4").Value = "SoC"
$ws.Range("C4").Value = "School of Computing"

# New rows 5-15, appended with additional faculties/schools
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "SCALE"
$ws.Range("C5").Value = "School of Continuing and Lifelong Education"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Dentistry"
$ws.Range("C6").Value = "Faculty of Dentistry"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "SDE"
$ws.Range("C7").Value = "School of Design & Environment"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "FoE"
$ws.Range("C8").Value = "Faculty of Engineering"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "ISEP"
$ws.Range("C9").Value = "Integrative Sciences and Engineering Programme"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Law"
$ws.Range("C10").Value = "Faculty of Law"

$ws.Range("A11").Value = 10
$ws.Range("C11").Value = "Yong Loo Lin School of Medicine"
$ws.Range("B11").Value = "Medicine"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Music"
$ws.Range("C12").Value = "Yong Siew Toh Conservatory of Music"

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Public Health"
$ws.Range("C13").Value = "Saw See Hock School of Public Health"

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Public Policy"
$ws.Range("C14").Value = "Lee Kuan Yew School of Public Policy"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Science"
$ws.Range("C15").Value = "Faculty of Science"

# Match the final selection state recorded in the workbook
$ws.Range("C12").Select()
